# Adds a new "2022-Q1" detail sheet (copied from the "2021-Q4" layout) right
# before the "总计" (total) summary sheet, fills it with the Q1-2022 fund
# holding data, and inserts a matching new row at the top of the "总计"
# sheet's data so the running summary includes the new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet by copying the "2021-Q4" sheet (same
#    column layout: 基金代码/基金名称/基金规模/股票总仓位/仓位占比/
#    持有市值(亿元)/仓位排名), placing the copy right before "总计".
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")
$template.Copy($total)
$newSheet = $wb.ActiveSheet
$newSheet.Name = "2022-Q1"

# NOTE: worksheet references returned by Item(...) are resolved to a fixed
# tab position at lookup time, not re-resolved live by name. Since the copy
# above inserted a sheet and shifted everything after it, re-fetch "总计"
# now that the sheet collection/order has changed, so later writes go to
# the right tab.
$total = $wb.Worksheets.Item("总计")

# Overwrite the copied data with the 2022-Q1 numbers. Columns B and D-G hold
# numeric-looking figures (fund codes, AUM, positions, ...) that the source
# sheet stores as plain TEXT (e.g. fund code "010783" keeps its leading
# zero) - a leading apostrophe forces text entry the same way typing
# '010783 into Excel does, instead of letting it auto-convert to a number.
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'010783"
$newSheet.Range("C2").Value = "德邦沪港深龙头混合A"
$newSheet.Range("D2").Value = "'0.93"
$newSheet.Range("E2").Value = "'81.58"
$newSheet.Range("F2").Value = "'3.27"
$newSheet.Range("G2").Value = "'0.0304"
$newSheet.Range("H2").Value = 8

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'010784"
$newSheet.Range("C3").Value = "德邦沪港深龙头混合C"
$newSheet.Range("D3").Value = "'0.27"
$newSheet.Range("E3").Value = "'81.58"
$newSheet.Range("F3").Value = "'3.27"
$newSheet.Range("G3").Value = "'0.0088"
$newSheet.Range("H3").Value = 8

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'161124"
$newSheet.Range("C4").Value = "易方达香港恒生综合小型股指数（QDII-LOF）A"
$newSheet.Range("D4").Value = "'0.28"
$newSheet.Range("E4").Value = "'92.62"
$newSheet.Range("F4").Value = "'1.91"
$newSheet.Range("G4").Value = "'0.0053"
$newSheet.Range("H4").Value = 3

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'006263"
$newSheet.Range("C5").Value = "易方达香港恒生综合小型股指数（QDII-LOF）C"
$newSheet.Range("D5").Value = "'0.06"
$newSheet.Range("E5").Value = "'92.62"
$newSheet.Range("F5").Value = "'1.91"
$newSheet.Range("G5").Value = "'0.0011"
$newSheet.Range("H5").Value = 3

# ---------------------------------------------------------------------
# 2) Insert a new row under the header of "总计" for the 2022-Q1 totals,
#    pushing the existing history rows down by one.
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.05

# Renumber the running index in column A for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
